$d = $word.ActiveDocument

# Avoid AutoFormat turning straight apostrophes into curly "smart" quotes
# while we edit the text below.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false
$word.AutoCorrect.ReplaceText = $false

function Replace-ExactText($searchText, $replaceText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false)
    if ($found) {
        $rng.Text = $replaceText
    }
    return $found
}

# 1. Identity document: Passeport -> Carte d'identité nationale (+ new number)
Replace-ExactText " Passeport  N°PP25342A   délivré le" " Carte d'identité nationale  N°AA-45467776-AQ   délivré le"

# 2. Issue date: 14 mars 2019 -> 12 juillet 2023
Replace-ExactText " 14 mars 2019  " " 12 juillet 2023  "

# 3. Issuing authority
Replace-ExactText " Direction générale de la documentation et l'immigation" " Forces nationales de police"

# 4. Loan purpose: Financer les etudes -> Mariage
Replace-ExactText "Financer les etudes" "Mariage"

# 5. Deposit guarantee wording update
Replace-ExactText " 150% de l'échéance" " 150 % du montant de l'échéance"

# 6. Remove the "Constitution de PEP / 20 000" bullet paragraph entirely.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Constitution de PEP*") {
        $p.Range.Delete()
        break
    }
}

# 7. Signature date: 04 décembre 2024 -> 11 décembre 2024
Replace-ExactText "04 décembre 2024" "11 décembre 2024"
